$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 483.33334
$ws.Range("I28").Value = 490.9091
$ws.Range("J28").Value = 400.0
$ws.Range("K28").Value = 490.9091
$ws.Range("L28").Value = 400.0
$ws.Range("M28").Value = -5.909100000000024
$ws.Range("N28").Value = -1370.0

$ws.Range("H107").Value = 736.6
$ws.Range("I107").Value = 613.2
$ws.Range("J107").Value = 860.0
$ws.Range("K107").Value = 613.2
$ws.Range("L107").Value = 860.0
$ws.Range("M107").Value = 1306.8
$ws.Range("N107").Value = -4700.0

$ws.Range("H129").Value = 8703.0
$ws.Range("I129").Value = 538.0
$ws.Range("K129").Value = 1614.0
$ws.Range("M129").Value = 3386.0

$ws.Range("H138").Value = 3097.5076
$ws.Range("I138").Value = 2166.1333
$ws.Range("K138").Value = 6498.3999
$ws.Range("M138").Value = -1358.3999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 816.7222
$ws.Range("I2").Value = 752.1177
$ws.Range("J2").Value = 874.5263
$ws.Range("K2").Value = 752.1177
$ws.Range("L2").Value = 874.5263
$ws.Range("M2").Value = -639.1177
$ws.Range("N2").Value = -1100.5263

$ws.Range("H45").Value = 1900.8422
$ws.Range("I45").Value = 2457.4285
$ws.Range("J45").Value = 1576.1666
$ws.Range("K45").Value = 2457.4285
$ws.Range("L45").Value = 1576.1666
$ws.Range("M45").Value = -2080.4285
$ws.Range("N45").Value = -2330.1666

$ws.Range("H105").Value = 0.0
$ws.Range("J105").Value = 0.0
$ws.Range("L105").Value = 0.0
$ws.Range("N105").Value = ""

$ws.Range("H116").Value = 816.7222
$ws.Range("I116").Value = 752.1177
$ws.Range("J116").Value = 874.5263
$ws.Range("K116").Value = 752.1177
$ws.Range("L116").Value = 874.5263
$ws.Range("M116").Value = 1541.8823
$ws.Range("N116").Value = -5462.5263

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 816.7222
$ws.Range("I3").Value = 752.1177
$ws.Range("J3").Value = 874.5263
$ws.Range("K3").Value = 752.1177
$ws.Range("L3").Value = 874.5263
$ws.Range("M3").Value = -638.1177
$ws.Range("N3").Value = -1102.5263

$ws.Range("H80").Value = 147.0
$ws.Range("I80").Value = 34.0
$ws.Range("J80").Value = 173.07692
$ws.Range("K80").Value = 34.0
$ws.Range("L80").Value = 173.07692
$ws.Range("M80").Value = 964.0
$ws.Range("N80").Value = -2169.07692

$ws.Range("H83").Value = 147.0
$ws.Range("I83").Value = 34.0
$ws.Range("J83").Value = 173.07692
$ws.Range("K83").Value = 170.0
$ws.Range("L83").Value = 865.3846
$ws.Range("M83").Value = 4822.0
$ws.Range("N83").Value = -10849.3846

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 35715188.0
$ws.Range("I16").Value = 62500628.0
$ws.Range("J16").Value = 1266.6666
$ws.Range("K16").Value = 62500628.0
$ws.Range("L16").Value = 1266.6666
$ws.Range("M16").Value = -62500341.0
$ws.Range("N16").Value = -1840.6666

$ws.Range("H68").Value = 0.0
$ws.Range("J68").Value = 0.0
$ws.Range("L68").Value = 0.0
$ws.Range("N68").Value = ""

$ws.Range("H71").Value = 0.0
$ws.Range("J71").Value = 0.0
$ws.Range("L71").Value = 0.0
$ws.Range("N71").Value = ""

$ws.Range("H74").Value = 0.0
$ws.Range("J74").Value = 0.0
$ws.Range("L74").Value = 0.0
$ws.Range("N74").Value = ""

$ws.Range("H75").Value = 0.0
$ws.Range("J75").Value = 0.0
$ws.Range("L75").Value = 0.0
$ws.Range("N75").Value = ""

$ws.Range("H77").Value = 0.0
$ws.Range("J77").Value = 0.0
$ws.Range("L77").Value = 0.0
$ws.Range("N77").Value = ""

$ws.Range("H78").Value = 0.0
$ws.Range("J78").Value = 0.0
$ws.Range("L78").Value = 0.0
$ws.Range("N78").Value = ""

$ws.Range("H87").Value = 0.0
$ws.Range("J87").Value = 0.0
$ws.Range("L87").Value = 0.0
$ws.Range("N87").Value = ""

$ws.Range("H90").Value = 0.0
$ws.Range("J90").Value = 0.0
$ws.Range("L90").Value = 0.0
$ws.Range("N90").Value = ""

$ws.Range("H92").Value = 0.0
$ws.Range("J92").Value = 0.0
$ws.Range("L92").Value = 0.0
$ws.Range("N92").Value = ""

$ws.Range("H95").Value = 0.0
$ws.Range("J95").Value = 0.0
$ws.Range("L95").Value = 0.0
$ws.Range("N95").Value = ""

$ws.Range("H108").Value = 0.0
$ws.Range("J108").Value = 0.0
$ws.Range("L108").Value = 0.0
$ws.Range("N108").Value = ""

$ws.Range("H113").Value = 35715188.0
$ws.Range("I113").Value = 62500628.0
$ws.Range("J113").Value = 1266.6666
$ws.Range("K113").Value = 62500628.0
$ws.Range("L113").Value = 1266.6666
$ws.Range("M113").Value = -62498458.0
$ws.Range("N113").Value = -5606.6666

$ws.Range("H119").Value = 1000.0
$ws.Range("I119").Value = 1000.0
$ws.Range("J119").Value = 0.0
$ws.Range("K119").Value = 1000.0
$ws.Range("L119").Value = 0.0
$ws.Range("M119").Value = 3838.0
$ws.Range("N119").Value = ""

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 485292.28
$ws.Range("I5").Value = 589.619
$ws.Range("K5").Value = 1768.857
$ws.Range("M5").Value = -1656.857

$ws.Range("H122").Value = 1382.0883
$ws.Range("I122").Value = 311.83334
$ws.Range("J122").Value = 2586.125
$ws.Range("K122").Value = 2806.50006
$ws.Range("L122").Value = 23275.125
$ws.Range("M122").Value = -356.5000600000003
$ws.Range("N122").Value = -28175.125

$ws.Range("H135").Value = 485292.28
$ws.Range("I135").Value = 589.619
$ws.Range("K135").Value = 5306.571
$ws.Range("M135").Value = -2771.571

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2780.5715
$ws.Range("I61").Value = 2320.2727
$ws.Range("J61").Value = 4468.3335
$ws.Range("K61").Value = 2320.2727
$ws.Range("L61").Value = 4468.3335
$ws.Range("M61").Value = -2118.2727
$ws.Range("N61").Value = -4872.3335

$ws.Range("H113").Value = 2780.5715
$ws.Range("I113").Value = 2320.2727
$ws.Range("J113").Value = 4468.3335
$ws.Range("K113").Value = 2320.2727
$ws.Range("L113").Value = 4468.3335
$ws.Range("M113").Value = -150.2727
$ws.Range("N113").Value = -8808.3335

$ws.Range("H141").Value = 0.0
$ws.Range("J141").Value = 0.0
$ws.Range("L141").Value = 0.0
$ws.Range("N141").Value = ""

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 8568.2
$ws.Range("J109").Value = 8874.75
$ws.Range("L109").Value = 8874.75
$ws.Range("N109").Value = -11648.75
